$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds a date serial number (24-Apr-2024 -> 24-May-2024)
$ws.Range("A1").Value = 45436

# Update the price column (D29:D32) with the new computed values
$ws.Range("D29").Value = 13023.612
$ws.Range("D30").Value = 15370.992
$ws.Range("D31").Value = 18326.952
$ws.Range("D32").Value = 21816.167
